$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "25.385.76"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -2.14%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.665.03"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -4.11%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9963"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.25%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "236.35"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -4.13%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9978"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.16%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4807"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -4.72%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2598"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -5.24%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06155"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.51%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07077"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.39%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.650.98"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.93%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "14.69"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -3.68%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.5869"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -10.54%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.365"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -8.85%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "74.33"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -3.80%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.9978"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.21%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.9974"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.09%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "25.350.67"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.27%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.000006692"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.12%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "11.42"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -4.16%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.862.52"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -4.95%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "4.376"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.92%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "8.604"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.46%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.316"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -3.22%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "134.56"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.38%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "15.12"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.30%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.385"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.89%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "104.81"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.94%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.675"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -6.98%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "3.967"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.50%  "

$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.07654"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -5.99%  "

$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.601"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.96%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.04355"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -8.26%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.9968"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.13%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.599"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.34%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.6025"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.31%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.9425"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -6.04%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.618"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.70%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.8467"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -3.86%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.9985"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.01504"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -6.73%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "99.05"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.70%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.815"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -7.62%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.3743"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -4.85%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "4.664"
$c.Style = "Normal"

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.1115"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -5.72%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "6.199"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -3.15%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.05248"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.59%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "29.47"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -4.50%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.217"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.84%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.08%  "
